$d = $word.ActiveDocument

# Locate the "Heading 1" paragraph (the Introduction heading that carries the
# "introduction" bookmark) by its text, then expand the hit to the whole
# paragraph (including its end-of-paragraph mark) so the replacement XML we
# splice in fully replaces that paragraph.
$hit = $d.Content
$found = $hit.Find.Execute("Heading 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Heading 1' paragraph to split"
}
$hit.Expand(4)  # wdParagraph

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14 = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# Split the paragraph in two, the way Word does when you put the cursor right
# before the heading text and press Ctrl+Enter: the bookmark stays on the
# (now style-less) leading paragraph together with the manual page break -
# its pPr/rPr spells out the Heading1 character formatting explicitly, since
# the paragraph itself no longer carries the Heading1 style - while the
# trailing paragraph keeps the Heading1 style and the heading text, now
# preceded by the rendered-page-break marker.
$xml = '<w:p ' + $w + '><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:szCs w:val="32"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="introduction"/><w:r><w:br w:type="page"/></w:r></w:p>' + `
       '<w:p ' + $w + ' ' + $w14 + '><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Heading 1</w:t></w:r></w:p>'

$hit.InsertXML($xml)
